$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that Word had left at the very start
#    of the document (it marks the author's last edit position and is not
#    meaningful content).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fix a typo in the "Facit" (answer key) section: the equation for
#    problem 7a should read  x = 1/13  (the denominator was missing a
#    leading "1"), not  x = 1/3 . Locate the fraction robustly: it's the
#    only OMath object, appearing after the "Facit" heading, whose math
#    text collapses to "13" (i.e. numerator "1" and denominator "3").
# ---------------------------------------------------------------------------
$facitRange = $d.Content
[void]$facitRange.Find.Execute("Facit", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
$facitEnd = $facitRange.End

$target = $null
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    # OMath run text picked up via Range.Text is interleaved with stray
    # carriage-return field/arg separators - strip them before comparing.
    $clean = $om.Range.Text.Replace([char]13, "")
    if ($om.Range.Start -ge $facitEnd -and $clean -eq "13") {
        $target = $om
    }
}

if ($target -ne $null) {
    # Rebuild the fraction as 1 / 13: keep the numerator "1" run untouched,
    # and split the denominator into two runs ("1" then "3") with the
    # "_GoBack" bookmark sitting right where the new digit was typed -
    # exactly between the two denominator runs.
    $xmlFrag = '<m:oMath>' + `
                 '<m:f>' + `
                   '<m:fPr>' + `
                     '<m:ctrlPr>' + `
                       '<w:rPr>' + `
                         '<w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>' + `
                         '<w:i/>' + `
                         '<w:sz w:val="32"/>' + `
                       '</w:rPr>' + `
                     '</m:ctrlPr>' + `
                   '</m:fPr>' + `
                   '<m:num>' + `
                     '<m:r>' + `
                       '<w:rPr>' + `
                         '<w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>' + `
                         '<w:sz w:val="32"/>' + `
                       '</w:rPr>' + `
                       '<m:t>1</m:t>' + `
                     '</m:r>' + `
                   '</m:num>' + `
                   '<m:den>' + `
                     '<m:r>' + `
                       '<w:rPr>' + `
                         '<w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>' + `
                         '<w:sz w:val="32"/>' + `
                       '</w:rPr>' + `
                       '<m:t>1</m:t>' + `
                     '</m:r>' + `
                     '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
                     '<w:bookmarkEnd w:id="0"/>' + `
                     '<m:r>' + `
                       '<w:rPr>' + `
                         '<w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>' + `
                         '<w:sz w:val="32"/>' + `
                       '</w:rPr>' + `
                       '<m:t>3</m:t>' + `
                     '</m:r>' + `
                   '</m:den>' + `
                 '</m:f>' + `
               '</m:oMath>'
    $target.Range.InsertXML($xmlFrag)
}
